$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This handback-status report regenerated with a new pair of e2e files:
#   0a34ea6d-c30a-4cf0-a4b2-c4e7a97d0f6b.md  -> 0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.md
#   e3817207-2a06-43a7-a18f-c31df4325e01.md  -> ffff2fd1e58b-2408-4bc1-9197-ecfdddf93ad2.md
# plus refreshed handoff/handback timestamps and a single shared xlf content
# hash (e24154a4794c8801e1eece2bf03a3390907cfa15) used for both locales.
# ---------------------------------------------------------------------------

# ---------------- Overview sheet ----------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.md"
$wsOverview.Range("B2").Value = "e2e\0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.md"
$wsOverview.Range("G2").Value = "2016-08-30 01:04:40"

$wsOverview.Range("A3").Value = "ffff2fd1e58b-2408-4bc1-9197-ecfdddf93ad2.md"
$wsOverview.Range("B3").Value = "e2e\ffff2fd1e58b-2408-4bc1-9197-ecfdddf93ad2.md"
$wsOverview.Range("G3").Value = "2016-08-30 01:04:40"

# Hyperlinks.Delete() on a worksheet clears every hyperlink on that sheet, so
# rebuild all of them (targets are unchanged - only the display text moves to
# the new file names).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e44a615a2b135b2edcef7dbb28afff6e6b091e62/e2e/0a34ea6d-c30a-4cf0-a4b2-c4e7a97d0f6b.md", "", "", "e2e\0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e44a615a2b135b2edcef7dbb28afff6e6b091e62/e2e/e3817207-2a06-43a7-a18f-c31df4325e01.md", "", "", "e2e\ffff2fd1e58b-2408-4bc1-9197-ecfdddf93ad2.md")

# ---------------- zh-cn sheet ----------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.md"
$wsZhCn.Range("G2").Value = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.e24154a4794c8801e1eece2bf03a3390907cfa15.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-30 01:04:35"
$wsZhCn.Range("I2").Value = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.md"
$wsZhCn.Range("J2").Value = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.e24154a4794c8801e1eece2bf03a3390907cfa15.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-30 01:04:57"

$wsZhCn.Range("A3").Value = "ffff2fd1e58b-2408-4bc1-9197-ecfdddf93ad2.md"
$wsZhCn.Range("G3").Value = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.e24154a4794c8801e1eece2bf03a3390907cfa15.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-30 01:04:35"
$wsZhCn.Range("I3").Value = "ffff2fd1e58b-2408-4bc1-9197-ecfdddf93ad2.md"
$wsZhCn.Range("J3").Value = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.e24154a4794c8801e1eece2bf03a3390907cfa15.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-30 01:04:57"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e44a615a2b135b2edcef7dbb28afff6e6b091e62/e2e/0a34ea6d-c30a-4cf0-a4b2-c4e7a97d0f6b.md", "", "", "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/99837fd668ffebf7726a8d233d9f29d188650913/e2e/0a34ea6d-c30a-4cf0-a4b2-c4e7a97d0f6b.md", "", "", "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e44a615a2b135b2edcef7dbb28afff6e6b091e62/e2e/e3817207-2a06-43a7-a18f-c31df4325e01.md", "", "", "ffff2fd1e58b-2408-4bc1-9197-ecfdddf93ad2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/99837fd668ffebf7726a8d233d9f29d188650913/e2e/e3817207-2a06-43a7-a18f-c31df4325e01.md", "", "", "ffff2fd1e58b-2408-4bc1-9197-ecfdddf93ad2.md")

# ---------------- de-de sheet ----------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.md"
$wsDeDe.Range("G2").Value = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.e24154a4794c8801e1eece2bf03a3390907cfa15.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-30 01:04:40"
$wsDeDe.Range("I2").Value = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.md"
$wsDeDe.Range("J2").Value = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.e24154a4794c8801e1eece2bf03a3390907cfa15.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-30 01:05:13"

$wsDeDe.Range("A3").Value = "ffff2fd1e58b-2408-4bc1-9197-ecfdddf93ad2.md"
$wsDeDe.Range("G3").Value = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.e24154a4794c8801e1eece2bf03a3390907cfa15.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-30 01:04:40"
$wsDeDe.Range("I3").Value = "ffff2fd1e58b-2408-4bc1-9197-ecfdddf93ad2.md"
$wsDeDe.Range("J3").Value = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.e24154a4794c8801e1eece2bf03a3390907cfa15.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-30 01:05:13"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e44a615a2b135b2edcef7dbb28afff6e6b091e62/e2e/0a34ea6d-c30a-4cf0-a4b2-c4e7a97d0f6b.md", "", "", "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7052f928cd6a096928a4a20ea91149cf67c23133/e2e/0a34ea6d-c30a-4cf0-a4b2-c4e7a97d0f6b.md", "", "", "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e44a615a2b135b2edcef7dbb28afff6e6b091e62/e2e/e3817207-2a06-43a7-a18f-c31df4325e01.md", "", "", "ffff2fd1e58b-2408-4bc1-9197-ecfdddf93ad2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7052f928cd6a096928a4a20ea91149cf67c23133/e2e/e3817207-2a06-43a7-a18f-c31df4325e01.md", "", "", "ffff2fd1e58b-2408-4bc1-9197-ecfdddf93ad2.md")
